# SQLAssignments.docx edit:
#   1) "How many products are sold by each suppliers" -- move "each " so that
#      it belongs to the run containing "suppliers" instead of the run
#      containing "...sold by ". The rendered text is unchanged, only the
#      run split point moves.
#   2) Merge the "Part" run and the following " " run (in "Part 2 : Normalize
#      the below data") into a single run "Part ".
#   3) Merge the "PET " run and the "HEALTH HISTORY " run (which currently
#      have a "_GoBack" bookmark sitting between them) into a single run
#      "PET HEALTH HISTORY ", dropping the now-pointless bookmark.

$d = $word.ActiveDocument

# Helper: forcing a Range.Text assignment to actually rewrite/merge the
# underlying run(s) even when the final text happens to equal the
# concatenation of the run(s) it replaces (a same-text assignment would
# otherwise be treated as a no-op and leave the run split / bookmark alone).
function Set-RangeTextForce($range, [string]$finalText) {
    $range.Text = "zzTEMPzz"
    $range.Text = $finalText
}

# ---------------------------------------------------------------------
# Edit 1: "...sold by each " / "suppliers"  ->  "...sold by " / "each suppliers"
# ---------------------------------------------------------------------
$rEach = $d.Content
$rEach.Find.Execute("each ") | Out-Null
$rEach.Delete()

$rSuppliers = $d.Content
$rSuppliers.Find.Execute("suppliers") | Out-Null
$rSuppliers.Text = "each suppliers"

# ---------------------------------------------------------------------
# Edit 2: "Part" + " " runs -> single "Part " run (heading "Part 2 : Normalize...")
# ---------------------------------------------------------------------
$rPartHeading = $d.Content
$rPartHeading.Find.Execute("Part 2 : Normalize") | Out-Null
$rPartSpace = $d.Range($rPartHeading.Start, $rPartHeading.Start + 5)
Set-RangeTextForce $rPartSpace "Part "

# ---------------------------------------------------------------------
# Edit 3: "PET " + bookmark(_GoBack) + "HEALTH HISTORY " -> single run,
# bookmark removed.
# ---------------------------------------------------------------------
$rPetHealth = $d.Content
$rPetHealth.Find.Execute("PET HEALTH HISTORY") | Out-Null
Set-RangeTextForce $rPetHealth "PET HEALTH HISTORY"

$d.Save()
